$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 55; $r++) {
    $ws.Cells.Item($r, 4).Value = "2024-07-06T13:10:00.000Z"
}
for ($r = 56; $r -le 115; $r++) {
    $ws.Cells.Item($r, 4).Value = "2024-07-06T13:11:00.000Z"
}
for ($r = 116; $r -le 121; $r++) {
    $ws.Cells.Item($r, 4).Value = "2024-07-06T13:12:00.000Z"
}
